# Re-fixed sending of player name to server.
# Fill in the two previously-blank Time Log entries (rows 53-54 on Sheet1)
# with the actual logged session data, and move the sheet's active
# selection to reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 53: 2014-09-19, 8:45 PM - 9:30 PM, 10 min interruption, Coding
$ws.Range("A53").Value = 41901
$ws.Range("B53").Value = 0.86458333333333337
$ws.Range("C53").Value = 0.89583333333333337
$ws.Range("D53").Value = 10
$ws.Range("F53").Value = "Coding"

# Row 54: 2014-09-19, 10:15 PM - 11:29 PM, 10 min interruption, Testing
$ws.Range("A54").Value = 41901
$ws.Range("B54").Value = 0.92708333333333337
$ws.Range("C54").Value = 0.9784722222222223
$ws.Range("D54").Value = 10
$ws.Range("F54").Value = "Testing"

# The E column holds a shared formula (=IF(...), (C-B)*24-D/60, "")) that
# was already present (evaluating to "" while the row was blank). Re-assign
# it to itself so the engine re-derives the cached result now that B/C/D
# are populated, instead of keeping the stale blank-branch value.
$ws.Range("E53").Formula = $ws.Range("E53").Formula
$ws.Range("E54").Formula = $ws.Range("E54").Formula

# Reflect the new scroll/selection position left behind by the edit.
$ws.Range("C55").Select()
